# edit.ps1 - applies the "arrows.pptx" activity update:
#  - slide 1: tiny text tweak (count(data= ...) -> count(x= ...))
#  - duplicate slide 1 -> new slide 2, rewritten for the 2021 dataset
#  - duplicate (old) slide 2 -> new slide 4, rewritten for the 2021 dataset,
#    plus one extra "names()" arrow shape
#  - the deck-wide "date updated" footer field bumps by one day on every
#    master/layout that carries it

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 (original): small text correction on the "count(...)" arrow
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$countShape1 = $slide1.Shapes.Item(3)
$countShape1.TextFrame.TextRange.Text = "count(x= sqf_2011, race) "

# ---------------------------------------------------------------------
# 2) Duplicate slide 1 -> becomes slide 2 (inserted right after slide 1).
#    Update its three arrows so they talk about sqf_2021 instead of
#    sqf_2011, with the new recode() mapping and SUSPECT_* columns.
# ---------------------------------------------------------------------
$slide1.Duplicate() | Out-Null
$slide2 = $p.Slides.Item(2)

# -- "select(...)" arrow --
$selShape = $slide2.Shapes.Item(1)
$selTr = $selShape.TextFrame.TextRange
$selTr.Paragraphs(1).Text = "select("
$selTr.Paragraphs(2).Text = "data = sqf_2021, "
$selTr.Paragraphs(3).Text = "SUSPECT_SEX, SUSPECT_RACE_DESCRIPTION, SUSPECT_REPORTED_AGE)"

# -- "recode(...)" arrow --
$recShape = $slide2.Shapes.Item(2)
$recTr = $recShape.TextFrame.TextRange

$recTr.Paragraphs(1).Text = "recode("
$recTr.Paragraphs(1).InsertAfter("`r" + "    sqf_2021`$race,") | Out-Null

$recTr.Paragraphs(3).Text = '    "BLACK HISPANIC" = "Hispanic",'
$recTr.Paragraphs(4).Text = '    "WHITE HISPANIC"= "Hispanic",'
$recTr.Paragraphs(5).Text = '    "BLACK"= "Black",'
$recTr.Paragraphs(6).Text = '    "WHITE"= "White",'
$recTr.Paragraphs(7).Text = '    "MIDDLE EASTERN/SOUTHWEST ASIAN"= "MESA",'
$recTr.Paragraphs(8).Text = '    "ASIAN / PACIFIC ISLANDER"= "API",'
$recTr.Paragraphs(9).Text = '    "AMERICAN INDIAN/ALASKAN NATIVE"= "AMAN",'
$recTr.Paragraphs(10).Text = '    "no data"= "no data",'
$recTr.Paragraphs(11).Text = '    .default = NA_character_'
$recTr.Paragraphs(11).InsertAfter("`r" + "  )") | Out-Null

# -- "count(...)" arrow --
$cntShape = $slide2.Shapes.Item(3)
$cntShape.TextFrame.TextRange.Text = "count(x= sqf_2021, race) "

# ---------------------------------------------------------------------
# 3) Duplicate the original second slide (now slide 3) -> becomes
#    slide 4. Update its arrows for sqf_2021 and append a new arrow
#    that renames the resulting columns.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$slide3.Duplicate() | Out-Null
$slide4 = $p.Slides.Item(4)

# -- "rename(...)" arrow --
$renShape = $slide4.Shapes.Item(1)
$renTr = $renShape.TextFrame.TextRange
$renTr.Paragraphs(1).Text = "rename("
$renTr.Paragraphs(2).Text = "data= sqf_2021, total = n )"

# -- "mutate(...)" arrow --
$mutShape = $slide4.Shapes.Item(2)
$mutTr = $mutShape.TextFrame.TextRange
$mutTr.Paragraphs(1).Text = "mutate("
$mutTr.Paragraphs(2).Text = "data = sqf_2021, percentage = total/sum(total) * 100 ) "

# -- "arrange(...)" arrow --
$arrShape = $slide4.Shapes.Item(3)
$arrShape.TextFrame.TextRange.Text = "arrange(data =sqf_2021 desc(percentage))"

# -- new "names(...)" arrow, cloned from the rename arrow so it keeps
#    the same shape style/fill as the rest of the deck --
$namesShape = $renShape.Duplicate()
$namesShape.Name = "Right Arrow 1"
$namesShape.Left = 6079138 / 12700.0
$namesShape.Top = 3487321 / 12700.0
$namesShape.Width = 4044251 / 12700.0
$namesShape.Height = 2914875 / 12700.0

$namesTr = $namesShape.TextFrame.TextRange
$namesTr.Paragraphs(1).Text = 'names(sqf_2021) <- c("sex","race","age") '
if ($namesTr.Paragraphs().Count -gt 1) {
    $namesTr.Paragraphs(2).Delete()
}

# ---------------------------------------------------------------------
# 4) The printed "date last touched" footer rolled from 10/10/23 to
#    10/11/23 on the slide master and on every slide layout.
# ---------------------------------------------------------------------
function Update-DateText($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "10/10/23") {
                $sh.TextFrame.TextRange.Text = "10/11/23"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateText($master)

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateText($layouts.Item($i))
}
